$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 114 (shifts old rows 114-148 down to 117-151)
$ws.Range("A114:A116").EntireRow.Insert()

# Fill in the 3 newly inserted rows with the "Artic Star" data point
# (columns A,B,C,E,F,G,H,I,J are constant across the whole sheet)
$commonA = 8
$commonB = "Terminal La Palmera de La Serena"
$commonC = "Coquimbo"
$commonE = 4
$commonF = "Fruta"
$commonG = 100103
$commonH = "Frutos de hueso (carozo)"
$commonI = 100103006
$commonJ = "Nectarín"

$newDate = 44543
$newR = "Región de O'Higgins"
$newQ = "`$/bins (420 kilos)"
$newT = 420

$rows = @(
    @{ Row = 114; L = "Especial"; M = 16; N = 435000; O = 440000; P = 437500; S = 1042 },
    @{ Row = 115; L = "Primera";  M = 20; N = 395000; O = 400000; P = 397500; S = 946 },
    @{ Row = 116; L = "Segunda";  M = 16; N = 355000; O = 360000; P = 357500; S = 851 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $commonA
    $ws.Range("B$row").Value = $commonB
    $ws.Range("C$row").Value = $commonC
    $ws.Range("D$row").Value = $newDate
    $ws.Range("E$row").Value = $commonE
    $ws.Range("F$row").Value = $commonF
    $ws.Range("G$row").Value = $commonG
    $ws.Range("H$row").Value = $commonH
    $ws.Range("I$row").Value = $commonI
    $ws.Range("J$row").Value = $commonJ
    $ws.Range("K$row").Value = "Artic Star"
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $newQ
    $ws.Range("R$row").Value = $newR
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $newT
}

# Ensure the date column keeps its date style (s="2") as in the rest of the sheet
$ws.Range("D114:D116").NumberFormat = $ws.Range("D117").NumberFormat
